$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.828054639203998
$ws.Range("C2").Value = 0.04792406147589645
$ws.Range("D2").Value = 0.5152381671841937
$ws.Range("E2").Value = 0.1639553252878763
$ws.Range("G2").Value = 0.002551623589767713
$ws.Range("I2").Value = 1.617277362494804
$ws.Range("J2").Value = 0.05827119115450818
$ws.Range("K2").Value = 0.5218678611569487
$ws.Range("L2").Value = 0.4534027135190968
$ws.Range("O2").Value = 6.916752145533223
$ws.Range("B3").Value = 0.7928720045402713
$ws.Range("C3").Value = 0.0455418629143125
$ws.Range("D3").Value = 0.5128221173465022
$ws.Range("E3").Value = 0.1640847198606945
$ws.Range("G3").Value = 0.002554523241422897
$ws.Range("I3").Value = 1.626806794587409
$ws.Range("J3").Value = 0.05851187378665479
$ws.Range("K3").Value = 0.488148608338463
$ws.Range("L3").Value = 0.4492969709061612
$ws.Range("O3").Value = 6.948975187960343
$ws.Range("B4").Value = 0.771585408027903
$ws.Range("C4").Value = 0.04406219162893876
$ws.Range("D4").Value = 0.5115524385223722
$ws.Range("E4").Value = 0.164220584608314
$ws.Range("G4").Value = 0.002556399658985489
$ws.Range("I4").Value = 1.633278702442084
$ws.Range("J4").Value = 0.05866908899268708
$ws.Range("K4").Value = 0.467616847792641
$ws.Range("L4").Value = 0.4469630400411688
$ws.Range("O4").Value = 6.971591475876039
$ws.Range("B5").Value = 0.7629910051239222
$ws.Range("C5").Value = 0.04345495228776031
$ws.Range("D5").Value = 0.5110889153889957
$ws.Range("E5").Value = 0.1642901684252642
$ws.Range("G5").Value = 0.002557188531364042
$ws.Range("I5").Value = 1.63607226028947
$ws.Range("J5").Value = 0.05873553328932513
$ws.Range("K5").Value = 0.4592937715077028
$ws.Range("L5").Value = 0.4460591041133739
$ws.Range("O5").Value = 6.981520162367445
$ws.Range("B6").Value = 0.761568766355424
$ws.Range("C6").Value = 0.04335386363100469
$ws.Range("D6").Value = 0.5110152056523276
$ws.Range("E6").Value = 0.1643025824169264
$ws.Range("G6").Value = 0.002557320987864036
$ws.Range("I6").Value = 1.636545565605353
$ws.Range("J6").Value = 0.05874671008423027
$ws.Range("K6").Value = 0.4579143911450956
$ws.Range("L6").Value = 0.4459118585386364
$ws.Range("O6").Value = 6.983211851166402
$ws.Range("B7").Value = 0.7714691758526442
$ws.Range("C7").Value = 0.04405401942457132
$ws.Range("D7").Value = 0.5115459689595951
$ws.Range("E7").Value = 0.1642214654303888
$ws.Range("G7").Value = 0.002556410199758831
$ws.Range("I7").Value = 1.633315744752842
$ws.Range("J7").Value = 0.05866997544918462
$ws.Range("K7").Value = 0.4675044218915616
$ws.Range("L7").Value = 0.4469506581179701
$ws.Range("O7").Value = 6.971722492706846
$ws.Range("B8").Value = 0.815858539375796
$ws.Range("C8").Value = 0.04710621223176048
$ws.Range("D8").Value = 0.5143608139802893
$ws.Range("E8").Value = 0.1639882500627579
$ws.Range("G8").Value = 0.002552603504314093
$ws.Range("I8").Value = 1.620434322368602
$ws.Range("J8").Value = 0.05835222346028335
$ws.Range("K8").Value = 0.5102060724505861
$ws.Range("L8").Value = 0.4519483173614134
$ws.Range("O8").Value = 6.927275288080068
$ws.Range("B9").Value = 0.9053861421480178
$ws.Range("C9").Value = 0.05295648838240652
$ws.Range("D9").Value = 0.5215723984467218
$ws.Range("E9").Value = 0.163977280847508
$ws.Range("G9").Value = 0.002545897161443933
$ws.Range("I9").Value = 1.600095677689495
$ws.Range("J9").Value = 0.05780373654097026
$ws.Range("K9").Value = 0.5952897995800015
$ws.Range("L9").Value = 0.4632276210750774
$ws.Range("O9").Value = 6.862562682809653
$ws.Range("B10").Value = 0.9726485448103119
$ws.Range("C10").Value = 0.05717240417305902
$ws.Range("D10").Value = 0.5278970306475372
$ws.Range("E10").Value = 0.1642398043820812
$ws.Range("G10").Value = 0.00254142783968889
$ws.Range("I10").Value = 1.588148505846171
$ws.Range("J10").Value = 0.05744592343195265
$ws.Range("K10").Value = 0.6586028526747612
$ws.Range("L10").Value = 0.4724108683556665
$ws.Range("O10").Value = 6.828686058940235
$ws.Range("B11").Value = 1.003565819817197
$ws.Range("C11").Value = 0.0590725036426818
$ws.Range("D11").Value = 0.5309961332822581
$ws.Range("E11").Value = 0.1644176753809425
$ws.Range("G11").Value = 0.002539493060426585
$ws.Range("I11").Value = 1.583362939669748
$ws.Range("J11").Value = 0.0572928806031916
$ws.Range("K11").Value = 0.6875762879836316
$ws.Range("L11").Value = 0.476782181737164
$ws.Range("O11").Value = 6.81623925243062
$ws.Range("B12").Value = 1.015318725752934
$ws.Range("C12").Value = 0.05978946476817271
$ws.Range("D12").Value = 0.5322015000085116
$ws.Range("E12").Value = 0.1644934063738077
$ws.Range("G12").Value = 0.002538774476537588
$ws.Range("I12").Value = 1.581644054259677
$ws.Range("J12").Value = 0.05723632079614305
$ws.Range("K12").Value = 0.6985720703800382
$ws.Range("L12").Value = 0.4784652403656651
$ws.Range("O12").Value = 6.811951849087393
$ws.Range("B13").Value = 1.012785527882471
$ws.Range("C13").Value = 0.05963516867662122
$ws.Range("D13").Value = 0.5319404901260043
$ws.Range("E13").Value = 0.164476724335028
$ws.Range("G13").Value = 0.002538928611605489
$ws.Range("I13").Value = 1.582010098163224
$ws.Range("J13").Value = 0.05724844003005458
$ws.Range("K13").Value = 0.6962028658267343
$ws.Range("L13").Value = 0.4781015320626807
$ws.Range("O13").Value = 6.812856278848102
$ws.Range("B14").Value = 1.004531836521352
$ws.Range("C14").Value = 0.05913153995800258
$ws.Range("D14").Value = 0.5310946628127198
$ws.Range("E14").Value = 0.1644237381349889
$ws.Range("G14").Value = 0.002539433660462575
$ws.Range("I14").Value = 1.583219656301637
$ws.Range("J14").Value = 0.05728819947885366
$ws.Range("K14").Value = 0.6884804364110266
$ws.Range("L14").Value = 0.4769200928465409
$ws.Range("O14").Value = 6.815877990370097
$ws.Range("B15").Value = 0.9994820783535658
$ws.Range("C15").Value = 0.05882271833124264
$ws.Range("D15").Value = 0.5305807074313407
$ws.Range("E15").Value = 0.1643923723663896
$ws.Range("G15").Value = 0.002539744848494321
$ws.Range("I15").Value = 1.583972694879073
$ws.Range("J15").Value = 0.05731273471497822
$ws.Range("K15").Value = 0.6837533565383467
$ws.Range("L15").Value = 0.4762000354857463
$ws.Range("O15").Value = 6.817784337053553
$ws.Range("B16").Value = 0.9706343923161
$ws.Range("C16").Value = 0.05704786998798284
$ws.Range("D16").Value = 0.527698953957028
$ws.Range("E16").Value = 0.1642293536789445
$ws.Range("G16").Value = 0.002541556255248763
$ws.Range("I16").Value = 1.588474313654224
$ws.Range("J16").Value = 0.05745612049301663
$ws.Range("K16").Value = 0.6567127835873521
$ws.Range("L16").Value = 0.4721290813474042
$ws.Range("O16").Value = 6.829559104401199
$ws.Range("B17").Value = 0.9530185410343677
$ws.Range("C17").Value = 0.05595450388226197
$ws.Range("D17").Value = 0.5259878653217953
$ws.Range("E17").Value = 0.1641442957658725
$ws.Range("G17").Value = 0.002542692634235122
$ws.Range("I17").Value = 1.591402160809892
$ws.Range("J17").Value = 0.05754657127695584
$ws.Range("K17").Value = 0.6401679220394385
$ws.Range("L17").Value = 0.4696812275430347
$ws.Range("O17").Value = 6.83754146446833
$ws.Range("B18").Value = 0.9429164625761643
$ws.Range("C18").Value = 0.05532395951135527
$ws.Range("D18").Value = 0.5250245947839431
$ws.Range("E18").Value = 0.1641008758455875
$ws.Range("G18").Value = 0.002543355508992973
$ws.Range("I18").Value = 1.593147296138298
$ws.Range("J18").Value = 0.05759951211423076
$ws.Range("K18").Value = 0.6306679791948113
$ws.Range("L18").Value = 0.4682915368035339
$ws.Range("O18").Value = 6.842411693098086
$ws.Range("B19").Value = 0.9395012624922572
$ws.Range("C19").Value = 0.05511018178043514
$ws.Range("D19").Value = 0.5247020415625059
$ws.Range("E19").Value = 0.164087120679909
$ws.Range("G19").Value = 0.002543581539692585
$ws.Range("I19").Value = 1.593748667344485
$ws.Range("J19").Value = 0.05761759442954251
$ws.Range("K19").Value = 0.6274542654643938
$ws.Range("L19").Value = 0.4678241498918112
$ws.Range("O19").Value = 6.844108593093665
$ws.Range("B20").Value = 0.9548906690574199
$ws.Range("C20").Value = 0.05607106741931034
$ws.Range("D20").Value = 0.5261678511720618
$ws.Range("E20").Value = 0.1641527809863028
$ws.Range("G20").Value = 0.002542570706826306
$ws.Range("I20").Value = 1.591084161722776
$ws.Range("J20").Value = 0.05753684787453484
$ws.Range("K20").Value = 0.6419274757709275
$ws.Range("L20").Value = 0.4699399179102386
$ws.Range("O20").Value = 6.836662856491898
$ws.Range("B21").Value = 1.0069549245095
$ws.Range("C21").Value = 0.05927953765298355
$ws.Range("D21").Value = 0.5313422405976667
$ws.Range("E21").Value = 0.1644390744038091
$ws.Range("G21").Value = 0.002539284933651159
$ws.Range("I21").Value = 1.582861847958235
$ws.Range("J21").Value = 0.05727648336430846
$ws.Range("K21").Value = 0.6907480496034282
$ws.Range("L21").Value = 0.477266358457527
$ws.Range("O21").Value = 6.8149788828689
$ws.Range("B22").Value = 1.041245088316458
$ws.Range("C22").Value = 0.06136149941815461
$ws.Range("D22").Value = 0.5349093251372778
$ws.Range("E22").Value = 0.1646749865525301
$ws.Range("G22").Value = 0.002537219500388543
$ws.Range("I22").Value = 1.578031914724015
$ws.Range("J22").Value = 0.05711444435685387
$ws.Range("K22").Value = 0.7227957455243086
$ws.Range("L22").Value = 0.4822162229209965
$ws.Range("O22").Value = 6.80328968993075
$ws.Range("B23").Value = 1.022919935578358
$ws.Range("C23").Value = 0.06025169060785629
$ws.Range("D23").Value = 0.5329885872592683
$ws.Range("E23").Value = 0.1645446200937997
$ws.Range("G23").Value = 0.002538314378967883
$ws.Range("I23").Value = 1.580560000221254
$ws.Range("J23").Value = 0.057200185801074
$ws.Range("K23").Value = 0.7056786051277868
$ws.Range("L23").Value = 0.4795596431951452
$ws.Range("O23").Value = 6.809301366693063
$ws.Range("B24").Value = 0.9540442004679619
$ws.Range("C24").Value = 0.05601837512820396
$ws.Range("D24").Value = 0.5260864158392309
$ws.Range("E24").Value = 0.164148927743426
$ws.Range("G24").Value = 0.002542625800455038
$ws.Range("I24").Value = 1.591227736358768
$ws.Range("J24").Value = 0.05754124089833201
$ws.Range("K24").Value = 0.6411319443629395
$ws.Range("L24").Value = 0.4698229091234509
$ws.Range("O24").Value = 6.83705919962577
$ws.Range("B25").Value = 0.8809036815543436
$ws.Range("C25").Value = 0.051388287871454
$ws.Range("D25").Value = 0.5194408866511111
$ws.Range("E25").Value = 0.1639326207283851
$ws.Range("G25").Value = 0.00254763067335293
$ws.Range("I25").Value = 1.605071409999546
$ws.Range("J25").Value = 0.05794416228581944
$ws.Range("K25").Value = 0.5721302262109589
$ws.Range("L25").Value = 0.4600184863321459
$ws.Range("O25").Value = 6.877667719547389
